# Hoàn thiện Ngoại Trú
# Update the test-case record (row 2) on both the "Data" and "Check" sheets
# with the new patient/case identifiers.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("A2").Value = 3009
$wsData.Range("E2").Value = 46200608009
$wsData.Range("X2").Value = "DN4127460130009"

$wsCheck = $wb.Worksheets.Item("Check")
$wsCheck.Range("A2").Value = 3009
$wsCheck.Range("C2").Value = "DN4127460130009"
